$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.04380589351058006
$ws.Range("C2").Value = 0.021202344447374344
$ws.Range("D2").Value = 0.012852024286985397
$ws.Range("E2").Value = 0.009358244016766548
$ws.Range("F2").Value = 0.0004067481495440006
$ws.Range("J2").Value = 0.1277914047241211
$ws.Range("K2").Value = 1.479499340057373
